# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (with fund holding data) right after
# the "总计" (summary) sheet, and adds a corresponding new row to the
# "总计" sheet summarizing the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q3 and shift the existing quarters down one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push the existing data rows (2..6) down to (3..7) - this carries the
# cell formatting/styles along with the values.
$summary.Range("A2:D6").Copy($summary.Range("A3:D7"))

# Write the new 2022-Q3 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.08

# Fix up the running index column (A) for the rows that moved down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计", before the
#    existing "2022-Q2" sheet.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Clone the layout/formatting of an existing same-shaped quarter sheet
# (header + 3 fund rows) so fonts/borders/column styling match exactly.
$template = $wb.Worksheets.Item("2022-Q1")
$template.Range("A1:H4").Copy($newSheet.Range("A1:H4"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 159851 华宝中证金融科技主题ETF
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159851"
$newSheet.Range("C2").Value = "华宝中证金融科技主题ETF"
$newSheet.Range("D2").Value = "'1.94"
$newSheet.Range("E2").Value = "'98.27"
$newSheet.Range("F2").Value = "'2.98"
$newSheet.Range("G2").Value = "'0.0578"
$newSheet.Range("H2").Value = 9

# Row 3 - 516100 华夏中证金融科技主题ETF
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'516100"
$newSheet.Range("C3").Value = "华夏中证金融科技主题ETF"
$newSheet.Range("D3").Value = "'0.51"
$newSheet.Range("E3").Value = "'96.79"
$newSheet.Range("F3").Value = "'2.96"
$newSheet.Range("G3").Value = "'0.0151"
$newSheet.Range("H3").Value = 8

# Row 4 - 516860 博时中证金融科技主题ETF
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'516860"
$newSheet.Range("C4").Value = "博时中证金融科技主题ETF"
$newSheet.Range("D4").Value = "'0.34"
$newSheet.Range("E4").Value = "'98.57"
$newSheet.Range("F4").Value = "'3.01"
$newSheet.Range("G4").Value = "'0.0102"
$newSheet.Range("H4").Value = 8

# Restore "总计" as the active sheet (matches the original workbook view).
$summary.Activate()
